$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1342.6
$ws.Range("I17").Value = 1049.5
$ws.Range("J17").Value = 1538
$ws.Range("K17").Value = 3148.5
$ws.Range("L17").Value = 4614
$ws.Range("M17").Value = -2980.5
$ws.Range("N17").Value = -4950
$ws.Range("H129").Value = 923.89233
$ws.Range("J129").Value = 956.78687
$ws.Range("L129").Value = 2870.36061
$ws.Range("N129").Value = -12870.36061
$ws.Range("H136").Value = 46985.973
$ws.Range("J136").Value = 46985.973
$ws.Range("L136").Value = 46985.973
$ws.Range("N136").Value = -57185.973

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10228.561
$ws.Range("I32").Value = 6251.68
$ws.Range("K32").Value = 6251.68
$ws.Range("M32").Value = -5964.68
$ws.Range("H92").Value = 29300
$ws.Range("J92").Value = 29300
$ws.Range("L92").Value = 29300
$ws.Range("N92").Value = -34292
$ws.Range("H109").Value = 26009.318
$ws.Range("J109").Value = 26009.318
$ws.Range("L109").Value = 26009.318
$ws.Range("N109").Value = -28783.318
$ws.Range("H122").Value = 2831.077
$ws.Range("I122").Value = 1709.4546
$ws.Range("K122").Value = 5128.3638
$ws.Range("M122").Value = -2678.3638
$ws.Range("H132").Value = 2434.2424
$ws.Range("I132").Value = 1119.6522
$ws.Range("J132").Value = 5457.8
$ws.Range("K132").Value = 3358.9566
$ws.Range("L132").Value = 16373.4
$ws.Range("M132").Value = -828.9566
$ws.Range("N132").Value = -21433.4
$ws.Range("H134").Value = 43985
$ws.Range("J134").Value = 43985
$ws.Range("L134").Value = 43985
$ws.Range("N134").Value = -54125
$ws.Range("H137").Value = 53040
$ws.Range("J137").Value = 53040
$ws.Range("L137").Value = 53040
$ws.Range("N137").Value = -63240

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H99").Value = 4492.647
$ws.Range("I99").Value = 2596.4285
$ws.Range("J99").Value = 5820
$ws.Range("K99").Value = 2596.4285
$ws.Range("L99").Value = 5820
$ws.Range("M99").Value = -1098.4285
$ws.Range("N99").Value = -8816
$ws.Range("H118").Value = 28890
$ws.Range("J118").Value = 28890
$ws.Range("L118").Value = 28890
$ws.Range("N118").Value = -32204
$ws.Range("H134").Value = 2858.66
$ws.Range("I134").Value = 1668.475
$ws.Range("J134").Value = 7619.4
$ws.Range("K134").Value = 5005.424999999999
$ws.Range("L134").Value = 22858.2
$ws.Range("M134").Value = -2470.424999999999
$ws.Range("N134").Value = -27928.2
$ws.Range("H137").Value = 33114.617
$ws.Range("J137").Value = 33114.617
$ws.Range("L137").Value = 33114.617
$ws.Range("N137").Value = -43314.617

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2124.9167
$ws.Range("I132").Value = 1000.4865
$ws.Range("K132").Value = 3001.4595
$ws.Range("M132").Value = -471.4594999999999
$ws.Range("H134").Value = 7859.579
$ws.Range("I134").Value = 10252.182
$ws.Range("J134").Value = 4569.75
$ws.Range("K134").Value = 30756.546
$ws.Range("L134").Value = 13709.25
$ws.Range("M134").Value = -28221.546
$ws.Range("N134").Value = -18779.25
$ws.Range("H137").Value = 50780
$ws.Range("J137").Value = 50780
$ws.Range("L137").Value = 50780
$ws.Range("N137").Value = -60980
$ws.Range("H139").Value = 47510
$ws.Range("J139").Value = 47510
$ws.Range("L139").Value = 47510
$ws.Range("N139").Value = -57790

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1163.0588
$ws.Range("I107").Value = 319.14285
$ws.Range("J107").Value = 5101.3335
$ws.Range("K107").Value = 957.4285500000001
$ws.Range("L107").Value = 15304.0005
$ws.Range("M107").Value = 962.5714499999999
$ws.Range("N107").Value = -19144.0005
$ws.Range("H113").Value = 822.8372000000001
$ws.Range("I113").Value = 713.82855
$ws.Range("J113").Value = 1299.75
$ws.Range("K113").Value = 2141.48565
$ws.Range("L113").Value = 3899.25
$ws.Range("M113").Value = 28.51435000000038
$ws.Range("N113").Value = -8239.25
$ws.Range("H122").Value = 2636.88
$ws.Range("I122").Value = 506.5238
$ws.Range("J122").Value = 3465.3518
$ws.Range("K122").Value = 4558.7142
$ws.Range("L122").Value = 31188.1662
$ws.Range("M122").Value = -2108.7142
$ws.Range("N122").Value = -36088.1662
$ws.Range("H131").Value = 8929473
$ws.Range("J131").Value = 895.7368
$ws.Range("L131").Value = 2687.2104
$ws.Range("N131").Value = -12767.2104
$ws.Range("H137").Value = 4806.25
$ws.Range("I137").Value = 6200
$ws.Range("J137").Value = 4607.143
$ws.Range("K137").Value = 18600
$ws.Range("L137").Value = 13821.429
$ws.Range("M137").Value = -13500
$ws.Range("N137").Value = -24021.429
$ws.Range("H139").Value = 3049.9033
$ws.Range("I139").Value = 1266.6875
$ws.Range("J139").Value = 4952
$ws.Range("K139").Value = 3800.0625
$ws.Range("L139").Value = 14856
$ws.Range("M139").Value = 1339.9375
$ws.Range("N139").Value = -25136
$ws.Range("H140").Value = 2288.3
$ws.Range("I140").Value = 1778.2307
$ws.Range("J140").Value = 2678.353
$ws.Range("K140").Value = 5334.6921
$ws.Range("L140").Value = 8035.059
$ws.Range("M140").Value = -154.6921000000002
$ws.Range("N140").Value = -18395.059

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 38500
$ws.Range("J48").Value = 38500
$ws.Range("L48").Value = 38500
$ws.Range("N48").Value = -39470
$ws.Range("H52").Value = 39500
$ws.Range("J52").Value = 39500
$ws.Range("L52").Value = 39500
$ws.Range("N52").Value = -40018
$ws.Range("H122").Value = 5198
$ws.Range("I122").Value = 4108.8887
$ws.Range("K122").Value = 12326.6661
$ws.Range("M122").Value = -9876.666100000002
$ws.Range("H123").Value = 15256.223
$ws.Range("J123").Value = 15256.223
$ws.Range("L123").Value = 15256.223
$ws.Range("N123").Value = -20156.223
$ws.Range("H137").Value = 53780
$ws.Range("J137").Value = 53780
$ws.Range("L137").Value = 53780
$ws.Range("N137").Value = -63980

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 287.42105
$ws.Range("I55").Value = 284.14285
$ws.Range("J55").Value = 296.6
$ws.Range("K55").Value = 284.14285
$ws.Range("L55").Value = 296.6
$ws.Range("M55").Value = -111.14285
$ws.Range("N55").Value = -642.6
$ws.Range("H75").Value = 42200
$ws.Range("J75").Value = 42200
$ws.Range("L75").Value = 42200
$ws.Range("N75").Value = -44072
$ws.Range("H78").Value = 42200
$ws.Range("J78").Value = 42200
$ws.Range("L78").Value = 126600
$ws.Range("N78").Value = -135960
$ws.Range("H122").Value = 3591.7026
$ws.Range("I122").Value = 3088.0293
$ws.Range("J122").Value = 9300
$ws.Range("K122").Value = 9264.0879
$ws.Range("L122").Value = 27900
$ws.Range("M122").Value = -6814.0879
$ws.Range("N122").Value = -32800
$ws.Range("H133").Value = 46722.152
$ws.Range("J133").Value = 46722.152
$ws.Range("L133").Value = 46722.152
$ws.Range("N133").Value = -51782.152

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5750276.5
$ws.Range("I132").Value = 3413.0286
$ws.Range("J132").Value = 14495504
$ws.Range("K132").Value = 10239.0858
$ws.Range("L132").Value = 43486512
$ws.Range("M132").Value = -7709.085800000001
$ws.Range("N132").Value = -43491572
$ws.Range("H135").Value = 48509.453
$ws.Range("J135").Value = 48509.453
$ws.Range("L135").Value = 48509.453
$ws.Range("N135").Value = -58649.453
